$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Advance the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Update price list values (step 1 and 2)
$ws.Range("D34").Value = 275
$ws.Range("D35").Value = 307.2
$ws.Range("D36").Value = 340.5
$ws.Range("D37").Value = 289
$ws.Range("D38").Value = 326.5
$ws.Range("D39").Value = 368.6
